# Update "想去人数" (F column) figures on the "展览" and "全部类型" sheets
# F3: 624 -> 623
# F4: 1564 -> 1565
# F5: 7330 -> 7333
# F7: 176 -> 177

$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)
    $ws.Range("F3").Value = 623
    $ws.Range("F4").Value = 1565
    $ws.Range("F5").Value = 7333
    $ws.Range("F7").Value = 177
}
